$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update Requirement text for the "대여소 상세 정보 조회" use case (row 3, column B):
# "...예약대기 신청, 문자 알림을 통해..." -> "...예약대기 신청하고 문자 알림을 통해..."
$ws.Range("B3").Value = "대여소 리스트에서 특정 대여소를 선택해 상세 정보 조회(대여소 이름, 대여소 위치, 사용 가능 자전거 목록 등), 대여소에 자전거가 없는 경우 예약대기 신청하고 문자 알림을 통해 해당 내용 수신"

# Update Requirement text for the "자전거 즉시 대여" use case (row 4, column B):
# "...즉시 대여, 문자 알림을 통해..." -> "...즉시 대여하고 문자 알림을 통해..."
$ws.Range("B4").Value = "대여소에 자전거가 남아 있는 경우 즉시 대여하고 문자 알림을 통해 해당 내용 수신"

# Update the active selection from C3 to C6
$ws.Range("C6").Select()
